$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose pasture type (column D) changes from "Sown Permanent Pasture" to "Natural Pasture"
$rowsToNatural = @(3, 5, 7, 9, 11, 13)
foreach ($r in $rowsToNatural) {
    $ws.Range("D$r").Value = "Natural Pasture"
}

# Rows whose pasture type (column D) changes from "Natural Pasture" to "Sown Permanent Pasture"
$rowsToSown = @(16, 18, 20, 22, 24, 26, 28, 30, 32, 34, 36, 38, 40, 42)
foreach ($r in $rowsToSown) {
    $ws.Range("D$r").Value = "Sown Permanent Pasture"
}

# Update the active selection to D7, matching the new cursor position recorded in the file
$ws.Range("D7").Select()
